$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "max" column (column C) entirely, shifting D and E left.
$ws.Range("C1").EntireColumn.Delete()

# Update the new B2 value (previously the "1" in the max column's pair cell).
$ws.Range("B2").Value = 3.213467280767539
